# Re-generate the "File" list (column A) and its index (row order),
# per the automatically-produced list/index regeneration described in
# the commit message. Column B ("Review date") values stay fixed to their
# row; only which file each row shows is updated/re-sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$files = @(
    'Routine_Care/Nursing for Arterial and Central Venous Lines.pdf',
    'Routine_Care/VTE_Prevention/TED Stocking Sizing.pdf',
    'Routine_Care/Faecal  incontinence skin care.pdf',
    'Breathing(Respiratory)/Equipment/IPPB using an ICU Ventilator.pdf',
    'Drugs/heparin_critical_care_only.pdf',
    'Routine_Care/Invasive Flush Systems.pdf',
    'GI_Liver_and_Transplant/Pancreatic Irrigation.pdf',
    'Neurological/Thiopentone levels.pdf',
    'Infection_and_sepsis/Ebola/Ebola.pdf',
    'GI_Liver_and_Transplant/Nasal bridle.pdf',
    'Breathing(Respiratory)/Equipment/AMBU AScope.pdf',
    'Cardiovascular/EZ-IO Intraosseus Access Device_pub_em.pdf',
    'Routine_Care/Central venous catheter removal.pdf',
    'Routine_Care/Tracheostomy_nursing_care.pdf',
    'Infection_and_sepsis/SARI/Suspected Influeza A Avian influenza H5N1 and SARS.pdf',
    'Infection_and_sepsis/SARI/MERS-CoV Guideline.pdf',
    'Infection_and_sepsis/SARI/Suspected Influeza A-H7N9 Guideline.pdf',
    'ECLS/Extra Corporeal Carbon Dioxide Removal.pdf',
    'Airway/Critical care extubation checklist.pdf',
    'Airway/Tracheostomy_Laryngectomy/Hospital_in-patients_with_a_Tracheostomy.pdf',
    'Trauma and Burns/Mangement of burns.pdf',
    'Drugs/diazepam_diazemuls.pdf',
    'End_of_life_care/Reasons to report a death to PF.pdf',
    'Breathing(Respiratory)/salbutamol and ipratroprium MDI.pdf',
    'Neurological/SOP -  Femoral site care.pdf',
    'Airway/Tracheostomy_Laryngectomy/Tracheostomy change in Critical Care.pdf',
    'Airway/Tracheostomy_Laryngectomy/Tracheostomy suctioning cleaning guideline.pdf',
    'Drugs/sodium_bicarbonate.pdf',
    'Post_op_care/Anticoagulation antiplatelet agents and epidural analgesia.pdf',
    'Post_op_care/Epidural top-up.pdf',
    'Airway/Emergency intubation checklist_em_pub.pdf',
    'Covid-19/WGH/CoVid intubation checklist WGH.pdf',
    'Covid-19/SJH/SJH COVID19 ITU Intubation Action Card.pdf',
    'Covid-19/SJH/SJH COVID19 ED Intubation Action Card.pdf',
    'Drugs/heparin for Haemofiltration.pdf',
    'Diabetes_and_Glucose/Hyperosmolar Hyperglycaemic State.pdf',
    'Drugs/fentanyl.pdf',
    'Airway/Tracheostomy_Laryngectomy/Tracheostomy guideline.pdf',
    'Covid-19/WGH/WGH_CT_Transfer_May.pdf',
    'Cardiovascular/GJNH Acute Heart Failure Referral Form.pdf',
    'Organ_donation/Donation after circulatory death.pdf',
    'Airway/Percutaneous tracheostomy checklist.pdf',
    'Delirium/Risk assessment posi mit.pdf',
    'Delirium/Managing a Potentially Violent Patient.pdf',
    'Infection_and_sepsis/SOP Ultrasound Cleaning.pdf',
    'GI_Liver_and_Transplant/Treatment of constipation.pdf',
    'Airway/McGrath Mac.pdf',
    'GI_Liver_and_Transplant/Abdominal pressure measurement.pdf',
    'Breathing(Respiratory)/HFNO.pdf',
    'Drugs/ketamine_in_asthma.pdf',
    'Delirium/Drugs Causing Delirium and Agitiation.pdf',
    'Neurological/Sub arachnoid haemorrhage management.pdf',
    'Airway/Anticipated difficult airway tool.pdf',
    'Airway/Tracheostomy_Laryngectomy/Tracheostomy safety box contents.pdf',
    'End_of_life_care/Documentation following death.pdf',
    'Drugs/zanamivir.pdf',
    'Routine_Care/bBraun Spaceplus Failure EMERGENCY ACTION CARD_em.pdf',
    'Drugs/insulin.pdf',
    'Breathing(Respiratory)/Equipment/HFNO Set Up.pdf',
    'Breathing(Respiratory)/Inhaled Nitrous Oxide.pdf',
    'Breathing(Respiratory)/Equipment/APRV.pdf',
    'Cardiovascular/Steroids for Septic Shock.pdf',
    'Post_op_care/Epidural Haematoma.pdf',
    'Neurological/SOP for review of Neurosurgical patients in ITU by neurosurgical team.pdf',
    'Breathing(Respiratory)/Equipment/T piece Y piece.pdf',
    'Policies_and_admin/General Critical Care Interaction with HEPMA_pub.pdf',
    'Drugs/midazolam and thiopental levels.pdf',
    'Breathing(Respiratory)/Equipment/HFNO through ventilator.pdf',
    'Routine_Care/VTE_Prevention/Dalteparin_thromboprophylaxis.pdf',
    'Post_op_care/Adult Scoliosis Spinal Surgery Post-Op Care.pdf',
    'Post_op_care/Post op care pharyngo-laryngo-oesphagectomy PLOG.pdf',
    'GI_Liver_and_Transplant/Nasogastric feeding protocol.pdf',
    'Diabetes_and_Glucose/Intravenous Insulin Therapy (not for DKA or HHS).pdf',
    'Drugs/Antibiotic doses in CVVHD.pdf',
    'GI_Liver_and_Transplant/Jejunostomy feeding protocol.pdf',
    'GI_Liver_and_Transplant/Nasojejunal feeding protocol.pdf',
    'ECLS/RIE ECLS Anti Xa Protocol.pdf',
    'Infection_and_sepsis/Winter Infections Stepdown Guidance.pdf',
    'Drugs/vasopressin_sepsis.pdf',
    'Drugs/vasopressin organ donation.pdf',
    'Transfer/ACCP Transfers.pdf',
    'Covid-19/videos/Donning and Doffing Video.pdf',
    'Breathing(Respiratory)/Equipment/Ventilators Circuits Filters and Closed Suction - Set up and Maintenance.pdf',
    'Drugs/piperacillin_tazobactam extended_infusion.pdf',
    'Breathing(Respiratory)/Equipment/Bipap V60.pdf',
    'Breathing(Respiratory)/CPAP.pdf',
    'Infection_and_sepsis/Infection indications for IVIG.pdf',
    'Covid-19/Covid 19 Death Certification Guideline.pdf',
    'Procedures/CVC Guidance/Securing CVCs.pdf',
    'Routine_Care/Video Communication.pdf',
    'Neurological/Treatment of status epilepticus.pdf',
    'Drugs/isoprenaline.pdf',
    'Cardiovascular/Cardiogenic Shock.pdf',
    'Haematology_CAR-T/Haem_ICU_transfer.pdf',
    'Infection_and_sepsis/SARI/Influenza A and B Virology Sampling and Oseltamivir Dose.pdf',
    'Cardiovascular/Management of hypertension within Critical Care.pdf',
    'Drugs/aminophylline.pdf',
    'Haematology_CAR-T/CRS.pdf',
    'Drugs/phenytoin.pdf',
    'Drugs/pancuronium.pdf',
    'Haematology_CAR-T/ICANS.pdf',
    'Drugs/rocuronium.pdf',
    'Drugs/nimodipine.pdf',
    'Drugs/milrinone.pdf',
    'Policies_and_admin/General Critical Care SOP_pub.pdf',
    'Ethics_and_Law/Care at the End of Life (FICM).pdf',
    'GI_Liver_and_Transplant/ICU - Upper GI bleeding (Endoscopy guideline).pdf',
    'Infection_and_sepsis/Trip Out of Unit infection guidance.pdf',
    'Ethics_and_Law/DNACPR policy for Scotland.pdf',
    'End_of_life_care/CMO & NRS Guidance for Doctors completing MCCD.pdf',
    'Covid-19/COVID 19 ICM guidance basic goals_June_2022.pdf',
    'Drugs/clonidine.pdf',
    'Organ_donation/Organ Retrieval SOP.pdf',
    'Neurological/Critical Care MRI Procedure_pub.pdf',
    'Neurological/Management of traumatic brain injury.pdf',
    'Airway/Cook Staged Extubation Set.pdf',
    'Drugs/noradrenaline (central).pdf',
    'Post_op_care/Epidural hypotension.pdf',
    'Drugs/dexmedetomidine.pdf',
    'Drugs/glyceryl_trinitrate.pdf',
    'GI_Liver_and_Transplant/Fulminant Liver Failure.pdf',
    'Breathing(Respiratory)/Equipment/Passy Muir Valve.pdf',
    'GI_Liver_and_Transplant/Confirmation of Nasogastric Tube Position.pdf',
    'Infection_and_sepsis/Antifungal guidance in critical care.pdf',
    'Transfer/Transfer Outdoors to Garden Guideline.pdf',
    'Cardiovascular/Intra Aortic Balloon Pump Guideline_pub.pdf',
    'Cardiovascular/Intra Aortic Balloon Pump Bedside Checks_pub.pdf',
    'Breathing(Respiratory)/ARDS Strategy.pdf',
    'Drugs/adrenaline.pdf',
    'Drugs/dobutamine.pdf',
    'Drugs/hydralazine.pdf',
    'Post_op_care/Major OMFS Free Flap.pdf',
    'Drugs/Alteplase for massive PE.pdf',
    'Drugs/alfentanil.pdf',
    'Drugs/magnesium.pdf',
    'Drugs/vancomycin.pdf',
    'Drugs/all IV drug infusion information.pdf',
    'Drugs/neostigmine.pdf',
    'Drugs/labetalol.pdf',
    'Neurological/Intrathecal policy RIE feb 2023.1.pdf',
    'Infection_and_sepsis/Initial investigation and management in unidentified Infections.pdf',
    'Cardiovascular/Management of Acute Type B Aortic Dissection Guideline.pdf',
    'Drugs/midazolam.pdf',
    'Procedures/CVC Guidance/CVC NHL  April 2023.pdf',
    'Drugs/potassium.pdf',
    'Drugs/salbutamol.pdf',
    'Routine_Care/ICU Eye Care Guideline.pdf',
    'Drugs/amiodarone.pdf',
    'Procedures/Arterial Line insertion for ACCPs.pdf',
    'Drugs/nicardipine.pdf',
    'Drugs/phenobarbitone.pdf',
    'Drugs/morphine.pdf',
    'Drugs/noradrenaline (peripheral).pdf',
    'Drugs/phenylephrine.pdf',
    'Neurological/Ventriculitis Guideline.pdf',
    'Drugs/Epoprostenol.pdf',
    'Breathing(Respiratory)/Manual Ventilation and MHI.pdf',
    'Cardiovascular/Cardiac Output Monitoring_pub .pdf',
    'Drugs/calcium.pdf',
    'Cardiovascular/Pulmonary_Embolism_and_DVT/Catheter directed thrombolysis of iliofemoral DVT alteplase_pub.pdf',
    'Drugs/Vancomycin Continuous Infusion Fluid Restricted.pdf',
    'Drugs/atracurium.pdf',
    'Airway/Tracheostomy_Laryngectomy/Decannulation Guidline.pdf',
    'Drugs/dalteparin_thromboprophylaxis.pdf',
    'Policies_and_admin/Pet Visitation.pdf',
    'Procedures/Inadvertent Catheter Placement Guideline.pdf',
    'End_of_life_care/Pet Visitation.pdf',
    'Drugs/valproate.pdf',
    'Drugs/stress ulcer prophylaxis.pdf',
    'Drugs/ketamine_for_status epilepticus.pdf',
    'Drugs/Phosphate.pdf',
    'Breathing(Respiratory)/Proning Guideline.pdf',
    'GI_Liver_and_Transplant/Prokinetics in ICU.pdf',
    'Drugs/Thiopentone.pdf',
    'Drugs/Prokinetics in ICU.pdf',
    'Drugs/Octreotide.pdf',
    'Post_op_care/Prevention and treatment of paraplegia after major aortic procedures.pdf',
    'Procedures/ACCPs acquiring initial CVC competencies.pdf',
    'Procedures/ACCP CVC placement following completion of initial competencies.pdf',
    'Delirium/Violence and Agression Guidlines.pdf',
    'Post_op_care/Care of the Transgender Patient.pdf',
    'GI_Liver_and_Transplant/Plasma exchange in Acute Liver Failure Guideline.pdf',
    'Breathing(Respiratory)/Equipment/NIV through Drager Vent Set up in Critical Care.pdf',
    'End_of_life_care/Guideline following Sudden Cardiac Death where death occurs in ICU.pdf',
    'Cardiovascular/Guideline following Sudden Cardiac Death where death occurs in ICU.pdf',
    'Transfer/Transfer Guidelines.pdf',
    'Breathing(Respiratory)/Equipment/NIV through Nihon Kohden  Setup.pdf',
    'Infection_and_sepsis/Influenza in Critical Care.pdf',
    'Policies_and_admin/Anticipated Post op flow surgical patients.pdf',
    'Policies_and_admin/Discharge Home from Critical Care .pdf',
    'Policies_and_admin/Repatriaiton Checklist for Critical Care.pdf'
)

for ($i = 0; $i -lt $files.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $files[$i]
}

